# Fruta / hortaliza, semanal
# Apply weekly update: rows 442-446 get new price-report values (as of 2021-09-09 / serial 44448),
# and 5 new rows (447-451) are inserted to preserve the previous values that used to sit in
# rows 442-446 (so the historical records are not lost). The row that used to be 447
# (Royal Gala) shifts down to become row 452 with its values unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 5 blank rows before the current row 447; this pushes the existing row 447
#    (Royal Gala / Primera) down to row 452, keeping its contents intact.
$ws.Range("A447:A451").EntireRow.Insert()

# 2) Update rows 442-446 in place with this week's new values.

# Row 442: Fuji royal / Especial
$ws.Cells.Item(442, 4).Value = 44448
$ws.Cells.Item(442, 11).Value = "Fuji royal"
$ws.Cells.Item(442, 12).Value = "Especial"
$ws.Cells.Item(442, 13).Value = 360
$ws.Cells.Item(442, 14).Value = 10000
$ws.Cells.Item(442, 15).Value = 10000
$ws.Cells.Item(442, 16).Value = 10000
$ws.Cells.Item(442, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(442, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(442, 19).Value = 667
$ws.Cells.Item(442, 20).Value = 15

# Row 443: Fuji royal / Primera
$ws.Cells.Item(443, 4).Value = 44448
$ws.Cells.Item(443, 12).Value = "Primera"
$ws.Cells.Item(443, 13).Value = 200
$ws.Cells.Item(443, 14).Value = 8000
$ws.Cells.Item(443, 15).Value = 8000
$ws.Cells.Item(443, 16).Value = 8000
$ws.Cells.Item(443, 19).Value = 533

# Row 444: Granny Smith / Especial
$ws.Cells.Item(444, 4).Value = 44448
$ws.Cells.Item(444, 11).Value = "Granny Smith"
$ws.Cells.Item(444, 12).Value = "Especial"
$ws.Cells.Item(444, 13).Value = 180
$ws.Cells.Item(444, 14).Value = 10000
$ws.Cells.Item(444, 15).Value = 10000
$ws.Cells.Item(444, 16).Value = 10000
$ws.Cells.Item(444, 19).Value = 667

# Row 445: Granny Smith / Primera
$ws.Cells.Item(445, 4).Value = 44448
$ws.Cells.Item(445, 13).Value = 230
$ws.Cells.Item(445, 14).Value = 8000
$ws.Cells.Item(445, 15).Value = 8000
$ws.Cells.Item(445, 16).Value = 8000
$ws.Cells.Item(445, 19).Value = 533

# Row 446: Granny Smith / Segunda
$ws.Cells.Item(446, 4).Value = 44448
$ws.Cells.Item(446, 11).Value = "Granny Smith"
$ws.Cells.Item(446, 12).Value = "Segunda"
$ws.Cells.Item(446, 13).Value = 150

# 3) Fill the newly inserted rows 447-451 with the values that previously occupied
#    rows 442-446 (preserving the historical data).

# Row 447 (= old row 442): Granny Smith / Calibre 120
$ws.Cells.Item(447, 1).Value = 5
$ws.Cells.Item(447, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(447, 3).Value = "Maule"
$ws.Cells.Item(447, 4).Value = 44238
$ws.Cells.Item(447, 5).Value = 7
$ws.Cells.Item(447, 6).Value = "Fruta"
$ws.Cells.Item(447, 7).Value = 100104
$ws.Cells.Item(447, 8).Value = "Frutos de pepita"
$ws.Cells.Item(447, 9).Value = 100104002
$ws.Cells.Item(447, 10).Value = "Manzana"
$ws.Cells.Item(447, 11).Value = "Granny Smith"
$ws.Cells.Item(447, 12).Value = "Calibre 120"
$ws.Cells.Item(447, 13).Value = 400
$ws.Cells.Item(447, 14).Value = 16000
$ws.Cells.Item(447, 15).Value = 16000
$ws.Cells.Item(447, 16).Value = 16000
$ws.Cells.Item(447, 17).Value = "$/caja 18 kilos embalada"
$ws.Cells.Item(447, 18).Value = "Provincia de Cachapoal"
$ws.Cells.Item(447, 19).Value = 889
$ws.Cells.Item(447, 20).Value = 18

# Row 448 (= old row 443): Fuji royal / Segunda
$ws.Cells.Item(448, 1).Value = 5
$ws.Cells.Item(448, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(448, 3).Value = "Maule"
$ws.Cells.Item(448, 4).Value = 44399
$ws.Cells.Item(448, 5).Value = 7
$ws.Cells.Item(448, 6).Value = "Fruta"
$ws.Cells.Item(448, 7).Value = 100104
$ws.Cells.Item(448, 8).Value = "Frutos de pepita"
$ws.Cells.Item(448, 9).Value = 100104002
$ws.Cells.Item(448, 10).Value = "Manzana"
$ws.Cells.Item(448, 11).Value = "Fuji royal"
$ws.Cells.Item(448, 12).Value = "Segunda"
$ws.Cells.Item(448, 13).Value = 260
$ws.Cells.Item(448, 14).Value = 7000
$ws.Cells.Item(448, 15).Value = 7000
$ws.Cells.Item(448, 16).Value = 7000
$ws.Cells.Item(448, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(448, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(448, 19).Value = 467
$ws.Cells.Item(448, 20).Value = 15

# Row 449 (= old row 444): Pink Lady / Primera
$ws.Cells.Item(449, 1).Value = 5
$ws.Cells.Item(449, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(449, 3).Value = "Maule"
$ws.Cells.Item(449, 4).Value = 44399
$ws.Cells.Item(449, 5).Value = 7
$ws.Cells.Item(449, 6).Value = "Fruta"
$ws.Cells.Item(449, 7).Value = 100104
$ws.Cells.Item(449, 8).Value = "Frutos de pepita"
$ws.Cells.Item(449, 9).Value = 100104002
$ws.Cells.Item(449, 10).Value = "Manzana"
$ws.Cells.Item(449, 11).Value = "Pink Lady"
$ws.Cells.Item(449, 12).Value = "Primera"
$ws.Cells.Item(449, 13).Value = 240
$ws.Cells.Item(449, 14).Value = 7000
$ws.Cells.Item(449, 15).Value = 7000
$ws.Cells.Item(449, 16).Value = 7000
$ws.Cells.Item(449, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(449, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(449, 19).Value = 467
$ws.Cells.Item(449, 20).Value = 15

# Row 450 (= old row 445): Granny Smith / Primera
$ws.Cells.Item(450, 1).Value = 5
$ws.Cells.Item(450, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(450, 3).Value = "Maule"
$ws.Cells.Item(450, 4).Value = 44400
$ws.Cells.Item(450, 5).Value = 7
$ws.Cells.Item(450, 6).Value = "Fruta"
$ws.Cells.Item(450, 7).Value = 100104
$ws.Cells.Item(450, 8).Value = "Frutos de pepita"
$ws.Cells.Item(450, 9).Value = 100104002
$ws.Cells.Item(450, 10).Value = "Manzana"
$ws.Cells.Item(450, 11).Value = "Granny Smith"
$ws.Cells.Item(450, 12).Value = "Primera"
$ws.Cells.Item(450, 13).Value = 360
$ws.Cells.Item(450, 14).Value = 6000
$ws.Cells.Item(450, 15).Value = 7000
$ws.Cells.Item(450, 16).Value = 6444
$ws.Cells.Item(450, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(450, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(450, 19).Value = 430
$ws.Cells.Item(450, 20).Value = 15

# Row 451 (= old row 446): Pink Lady / Primera
$ws.Cells.Item(451, 1).Value = 5
$ws.Cells.Item(451, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(451, 3).Value = "Maule"
$ws.Cells.Item(451, 4).Value = 44400
$ws.Cells.Item(451, 5).Value = 7
$ws.Cells.Item(451, 6).Value = "Fruta"
$ws.Cells.Item(451, 7).Value = 100104
$ws.Cells.Item(451, 8).Value = "Frutos de pepita"
$ws.Cells.Item(451, 9).Value = 100104002
$ws.Cells.Item(451, 10).Value = "Manzana"
$ws.Cells.Item(451, 11).Value = "Pink Lady"
$ws.Cells.Item(451, 12).Value = "Primera"
$ws.Cells.Item(451, 13).Value = 130
$ws.Cells.Item(451, 14).Value = 6000
$ws.Cells.Item(451, 15).Value = 6000
$ws.Cells.Item(451, 16).Value = 6000
$ws.Cells.Item(451, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(451, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(451, 19).Value = 400
$ws.Cells.Item(451, 20).Value = 15

# Row 452 already holds the former row-447 data (Royal Gala / Primera), shifted
# down intact by the Insert() above - no further changes needed there.
